$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows at 268-272 (shifts old 268-278 down to 273-283)
$ws.Range("A268:A272").EntireRow.Insert()

# Row 268
$ws.Range("A268").Value = 8
$ws.Range("B268").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C268").Value = 'Coquimbo'
$ws.Range("D268").Value = 44578
$ws.Range("E268").Value = 4
$ws.Range("F268").Value = 'Fruta'
$ws.Range("G268").Value = 100103
$ws.Range("H268").Value = 'Frutos de hueso (carozo)'
$ws.Range("I268").Value = 100103001
$ws.Range("J268").Value = 'Cereza'
$ws.Range("K268").Value = 'Lapins'
$ws.Range("L268").Value = 'Primera'
$ws.Range("M268").Value = 360
$ws.Range("N268").Value = 8500
$ws.Range("O268").Value = 9000
$ws.Range("P268").Value = 8750
$ws.Range("Q268").Value = '$/bandeja 10 kilos'
$ws.Range("R268").Value = 'Provincia de Curicó'
$ws.Range("S268").Value = 875
$ws.Range("T268").Value = 10

# Row 269
$ws.Range("A269").Value = 8
$ws.Range("B269").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C269").Value = 'Coquimbo'
$ws.Range("D269").Value = 44578
$ws.Range("E269").Value = 4
$ws.Range("F269").Value = 'Fruta'
$ws.Range("G269").Value = 100103
$ws.Range("H269").Value = 'Frutos de hueso (carozo)'
$ws.Range("I269").Value = 100103001
$ws.Range("J269").Value = 'Cereza'
$ws.Range("K269").Value = 'Lapins'
$ws.Range("L269").Value = 'Segunda'
$ws.Range("M269").Value = 300
$ws.Range("N269").Value = 6500
$ws.Range("O269").Value = 7000
$ws.Range("P269").Value = 6750
$ws.Range("Q269").Value = '$/bandeja 10 kilos'
$ws.Range("R269").Value = 'Provincia de Curicó'
$ws.Range("S269").Value = 675
$ws.Range("T269").Value = 10

# Row 270
$ws.Range("A270").Value = 8
$ws.Range("B270").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C270").Value = 'Coquimbo'
$ws.Range("D270").Value = 44578
$ws.Range("E270").Value = 4
$ws.Range("F270").Value = 'Fruta'
$ws.Range("G270").Value = 100103
$ws.Range("H270").Value = 'Frutos de hueso (carozo)'
$ws.Range("I270").Value = 100103001
$ws.Range("J270").Value = 'Cereza'
$ws.Range("K270").Value = 'Santina'
$ws.Range("L270").Value = 'Primera'
$ws.Range("M270").Value = 400
$ws.Range("N270").Value = 8500
$ws.Range("O270").Value = 9000
$ws.Range("P270").Value = 8750
$ws.Range("Q270").Value = '$/bandeja 10 kilos'
$ws.Range("R270").Value = 'Provincia de Curicó'
$ws.Range("S270").Value = 875
$ws.Range("T270").Value = 10

# Row 271
$ws.Range("A271").Value = 8
$ws.Range("B271").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C271").Value = 'Coquimbo'
$ws.Range("D271").Value = 44578
$ws.Range("E271").Value = 4
$ws.Range("F271").Value = 'Fruta'
$ws.Range("G271").Value = 100103
$ws.Range("H271").Value = 'Frutos de hueso (carozo)'
$ws.Range("I271").Value = 100103001
$ws.Range("J271").Value = 'Cereza'
$ws.Range("K271").Value = 'Santina'
$ws.Range("L271").Value = 'Segunda'
$ws.Range("M271").Value = 400
$ws.Range("N271").Value = 6500
$ws.Range("O271").Value = 7000
$ws.Range("P271").Value = 6750
$ws.Range("Q271").Value = '$/bandeja 10 kilos'
$ws.Range("R271").Value = 'Provincia de Curicó'
$ws.Range("S271").Value = 675
$ws.Range("T271").Value = 10

# Row 272
$ws.Range("A272").Value = 8
$ws.Range("B272").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C272").Value = 'Coquimbo'
$ws.Range("D272").Value = 44578
$ws.Range("E272").Value = 4
$ws.Range("F272").Value = 'Fruta'
$ws.Range("G272").Value = 100103
$ws.Range("H272").Value = 'Frutos de hueso (carozo)'
$ws.Range("I272").Value = 100103001
$ws.Range("J272").Value = 'Cereza'
$ws.Range("K272").Value = 'Sweet Heart'
$ws.Range("L272").Value = 'Primera'
$ws.Range("M272").Value = 360
$ws.Range("N272").Value = 8500
$ws.Range("O272").Value = 9000
$ws.Range("P272").Value = 8750
$ws.Range("Q272").Value = '$/bandeja 10 kilos'
$ws.Range("R272").Value = 'Provincia de Curicó'
$ws.Range("S272").Value = 875
$ws.Range("T272").Value = 10
